$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 157.45454
$ws.Range("I33").Value = 86
$ws.Range("K33").Value = 86
$ws.Range("M33").Value = 143

$ws.Range("H80").Value = 844.6875
$ws.Range("I80").Value = 1548.5
$ws.Range("K80").Value = 4645.5
$ws.Range("M80").Value = -3647.5

$ws.Range("H83").Value = 844.6875
$ws.Range("I83").Value = 1548.5
$ws.Range("K83").Value = 13936.5
$ws.Range("M83").Value = -8944.5

$ws.Range("H113").Value = 3375.5
$ws.Range("I113").Value = 2751
$ws.Range("K113").Value = 2751
$ws.Range("M113").Value = 503

$ws.Range("H121").Value = 1819.1666
$ws.Range("J121").Value = 1883
$ws.Range("L121").Value = 5649
$ws.Range("N121").Value = -9143

$ws.Range("H132").Value = 11116368
$ws.Range("I132").Value = 14498532
$ws.Range("J132").Value = 3544.5715
$ws.Range("K132").Value = 43495596
$ws.Range("L132").Value = 10633.7145
$ws.Range("M132").Value = -43493066
$ws.Range("N132").Value = -15693.7145

$ws.Range("H137").Value = 1217.2
$ws.Range("I137").Value = 853.5238000000001
$ws.Range("K137").Value = 2560.5714
$ws.Range("M137").Value = -10.57140000000027

$ws.Range("H138").Value = 1551.2
$ws.Range("I138").Value = 817.8125
$ws.Range("J138").Value = 1923.7142
$ws.Range("K138").Value = 2453.4375
$ws.Range("L138").Value = 5771.142599999999
$ws.Range("M138").Value = 2686.5625
$ws.Range("N138").Value = -16051.1426

$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 184.33333
$ws.Range("I5").Value = 184.33333
$ws.Range("K5").Value = 184.33333
$ws.Range("M5").Value = -72.33332999999999

$ws.Range("H32").Value = 2741.442
$ws.Range("I32").Value = 2634.3157
$ws.Range("K32").Value = 2634.3157
$ws.Range("M32").Value = -2347.3157

$ws.Range("H61").Value = 1206.7778
$ws.Range("I61").Value = 980.1429000000001
$ws.Range("K61").Value = 980.1429000000001
$ws.Range("M61").Value = -768.1429000000001

$ws.Range("H74").Value = 837.46155
$ws.Range("I74").Value = 657.25
$ws.Range("K74").Value = 657.25
$ws.Range("M74").Value = 216.75

$ws.Range("H77").Value = 837.46155
$ws.Range("I77").Value = 657.25
$ws.Range("K77").Value = 3286.25
$ws.Range("M77").Value = 1081.75

$ws.Range("H132").Value = 2275.3438
$ws.Range("I132").Value = 1954.7693
$ws.Range("J132").Value = 3664.5
$ws.Range("K132").Value = 5864.3079
$ws.Range("L132").Value = 10993.5
$ws.Range("M132").Value = -3334.3079
$ws.Range("N132").Value = -16053.5

$ws.Range("H136").Value = 1206.7778
$ws.Range("I136").Value = 980.1429000000001
$ws.Range("K136").Value = 2940.4287
$ws.Range("M136").Value = -390.4287000000004

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 184.33333
$ws.Range("I4").Value = 184.33333
$ws.Range("K4").Value = 184.33333
$ws.Range("M4").Value = -69.33332999999999

$ws.Range("H20").Value = 1473.4286
$ws.Range("I20").Value = 1028.5
$ws.Range("K20").Value = 1028.5
$ws.Range("M20").Value = -781.5

$ws.Range("H86").Value = 3777.0322
$ws.Range("I86").Value = 3961.3044
$ws.Range("J86").Value = 3247.25
$ws.Range("K86").Value = 3961.3044
$ws.Range("L86").Value = 3247.25
$ws.Range("M86").Value = -2838.3044
$ws.Range("N86").Value = -5493.25

$ws.Range("H89").Value = 3777.0322
$ws.Range("I89").Value = 3961.3044
$ws.Range("J89").Value = 3247.25
$ws.Range("K89").Value = 19806.522
$ws.Range("L89").Value = 16236.25
$ws.Range("M89").Value = -14190.522
$ws.Range("N89").Value = -27468.25

$ws.Range("H134").Value = 7310.2104
$ws.Range("I134").Value = 1062.8572
$ws.Range("J134").Value = 24802.8
$ws.Range("K134").Value = 3188.5716
$ws.Range("L134").Value = 74408.39999999999
$ws.Range("M134").Value = -653.5715999999998
$ws.Range("N134").Value = -79478.39999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 62
$ws.Range("I7").Value = 10
$ws.Range("K7").Value = 10
$ws.Range("M7").Value = 103

$ws.Range("H31").Value = 1365.72
$ws.Range("I31").Value = 1318.4584
$ws.Range("K31").Value = 1318.4584
$ws.Range("M31").Value = -1023.4584

$ws.Range("H34").Value = 1365.72
$ws.Range("I34").Value = 1318.4584
$ws.Range("K34").Value = 1318.4584
$ws.Range("M34").Value = -1116.4584

$ws.Range("H58").Value = 1398.7142
$ws.Range("I58").Value = 1206.1818
$ws.Range("K58").Value = 1206.1818
$ws.Range("M58").Value = -1003.1818

$ws.Range("H99").Value = 1503.7333
$ws.Range("I99").Value = 1655.375
$ws.Range("J99").Value = 1330.4286
$ws.Range("K99").Value = 1655.375
$ws.Range("L99").Value = 1330.4286
$ws.Range("M99").Value = -157.375
$ws.Range("N99").Value = -4326.4286

$ws.Range("H107").Value = 394.6111
$ws.Range("I107").Value = 266.3846
$ws.Range("K107").Value = 266.3846
$ws.Range("M107").Value = 1653.6154

$ws.Range("H109").Value = 26285.857
$ws.Range("J109").Value = 26285.857
$ws.Range("L109").Value = 26285.857
$ws.Range("N109").Value = -28365.857

$ws.Range("H126").Value = 1503.7333
$ws.Range("I126").Value = 1655.375
$ws.Range("J126").Value = 1330.4286
$ws.Range("K126").Value = 4966.125
$ws.Range("L126").Value = 3991.2858
$ws.Range("M126").Value = -2496.125
$ws.Range("N126").Value = -8931.2858

$ws.Range("H132").Value = 2457.762
$ws.Range("I132").Value = 1758.2142
$ws.Range("K132").Value = 5274.642599999999
$ws.Range("M132").Value = -2744.642599999999

$ws.Range("H134").Value = 1256.1538
$ws.Range("I134").Value = 1131.6666
$ws.Range("K134").Value = 3394.9998
$ws.Range("M134").Value = -859.9998000000001

$ws.Range("H136").Value = 1398.7142
$ws.Range("I136").Value = 1206.1818
$ws.Range("K136").Value = 3618.5454
$ws.Range("M136").Value = -1068.5454

$ws.Range("H141").Value = 27995
$ws.Range("J141").Value = 27995
$ws.Range("L141").Value = 27995
$ws.Range("N141").Value = -38355

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 3431.1428
$ws.Range("J100").Value = 3431.1428
$ws.Range("L100").Value = 10293.4284
$ws.Range("N100").Value = -11915.4284

$ws.Range("H106").Value = 2571.6
$ws.Range("J106").Value = 2571.6
$ws.Range("L106").Value = 7714.799999999999
$ws.Range("N106").Value = -9606.799999999999

$ws.Range("H131").Value = 13335611
$ws.Range("I131").Value = 76923416
$ws.Range("J131").Value = 2684.8708
$ws.Range("K131").Value = 230770248
$ws.Range("L131").Value = 8054.6124
$ws.Range("M131").Value = -230765208
$ws.Range("N131").Value = -18134.6124

$ws.Range("H140").Value = 1942.4
$ws.Range("I140").Value = 1490
$ws.Range("J140").Value = 2998
$ws.Range("K140").Value = 4470
$ws.Range("L140").Value = 8994
$ws.Range("M140").Value = 710
$ws.Range("N140").Value = -19354

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2660
$ws.Range("J80").Value = 4350
$ws.Range("L80").Value = 4350
$ws.Range("N80").Value = -6346

$ws.Range("H83").Value = 2660
$ws.Range("J83").Value = 4350
$ws.Range("L83").Value = 21750
$ws.Range("N83").Value = -31734

$ws.Range("H97").Value = 743.8570999999999
$ws.Range("I97").Value = 743.8570999999999
$ws.Range("K97").Value = 743.8570999999999
$ws.Range("M97").Value = -247.8570999999999

$ws.Range("H132").Value = 1779.2188
$ws.Range("I132").Value = 1502.5264
$ws.Range("J132").Value = 2183.6155
$ws.Range("K132").Value = 4507.5792
$ws.Range("L132").Value = 6550.8465
$ws.Range("M132").Value = -1977.5792
$ws.Range("N132").Value = -11610.8465

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1811.5555
$ws.Range("I68").Value = 1567.6666
$ws.Range("J68").Value = 2299.3333
$ws.Range("K68").Value = 1567.6666
$ws.Range("L68").Value = 2299.3333
$ws.Range("M68").Value = -818.6666
$ws.Range("N68").Value = -3797.3333

$ws.Range("H71").Value = 1811.5555
$ws.Range("I71").Value = 1567.6666
$ws.Range("J71").Value = 2299.3333
$ws.Range("K71").Value = 7838.333000000001
$ws.Range("L71").Value = 11496.6665
$ws.Range("M71").Value = -4094.333000000001
$ws.Range("N71").Value = -18984.6665

$ws.Range("H82").Value = 1958
$ws.Range("I82").Value = 1896.6666
$ws.Range("K82").Value = 1896.6666
$ws.Range("M82").Value = -1535.6666

$ws.Range("H85").Value = 1958
$ws.Range("I85").Value = 1896.6666
$ws.Range("K85").Value = 1896.6666
$ws.Range("M85").Value = -648.6666

$ws.Range("H93").Value = 718
$ws.Range("I93").Value = 704.6
$ws.Range("J93").Value = 785
$ws.Range("K93").Value = 704.6
$ws.Range("L93").Value = 785
$ws.Range("M93").Value = 543.4
$ws.Range("N93").Value = -3281

$ws.Range("H132").Value = 30293.314
$ws.Range("I132").Value = 1082
$ws.Range("J132").Value = 64981.75
$ws.Range("K132").Value = 3246
$ws.Range("L132").Value = 194945.25
$ws.Range("M132").Value = -716
$ws.Range("N132").Value = -200005.25

$ws.Range("H136").Value = 948.6
$ws.Range("I136").Value = 945.8946999999999
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 2837.6841
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -287.6840999999999
$ws.Range("N136").Value = -8100
